$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Transport -> Rest, 800 -> 2400, date 45961... -> 45964...
$ws.Range("A2").Value = "Rest"
$ws.Range("B2").Value = 2400
$ws.Range("C2").Value = 45964.083333333336

# Remove row 3 entirely (Category "rent" row)
$ws.Range("A3:C3").Delete()
